$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for the rows involved in the cyclic swap
$row4  = @($ws.Range("A4").Value2,  $ws.Range("B4").Value2,  $ws.Range("C4").Value2)
$row5  = @($ws.Range("A5").Value2,  $ws.Range("B5").Value2,  $ws.Range("C5").Value2)
$row7  = @($ws.Range("A7").Value2,  $ws.Range("B7").Value2,  $ws.Range("C7").Value2)
$row14 = @($ws.Range("A14").Value2, $ws.Range("B14").Value2, $ws.Range("C14").Value2)
$row15 = @($ws.Range("A15").Value2, $ws.Range("B15").Value2, $ws.Range("C15").Value2)

# Apply cyclic permutation: 4 -> 14 -> 7 -> 15 -> 5 -> 4
$ws.Range("A14").Value2 = $row4[0]
$ws.Range("B14").Value2 = $row4[1]
$ws.Range("C14").Value2 = $row4[2]

$ws.Range("A7").Value2 = $row14[0]
$ws.Range("B7").Value2 = $row14[1]
$ws.Range("C7").Value2 = $row14[2]

$ws.Range("A15").Value2 = $row7[0]
$ws.Range("B15").Value2 = $row7[1]
$ws.Range("C15").Value2 = $row7[2]

$ws.Range("A5").Value2 = $row15[0]
$ws.Range("B5").Value2 = $row15[1]
$ws.Range("C5").Value2 = $row15[2]

$ws.Range("A4").Value2 = $row5[0]
$ws.Range("B4").Value2 = $row5[1]
$ws.Range("C4").Value2 = $row5[2]
